$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = '71.230.73'
$ws.Range("E2").Value = '  +2.46%  '

# Row 3
$ws.Range("D3").Value = '4.003.92'
$ws.Range("E3").Value = '  +1.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '529.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.51%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.61%  '

# Row 7
$ws.Range("E7").Value = '  -0.57%  '

# Row 8
$ws.Range("E8").Value = '  +0.09%  '

# Row 9
$ws.Range("E9").Value = '  +0.15%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.176'
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = '  -2.24%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.73'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.21%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.48%  '

# Row 14
$ws.Range("D14").Value = '4.639.67'
$ws.Range("E14").Value = '  +1.53%  '

# Row 15
$ws.Range("D15").Value = '3.999.06'
$ws.Range("E15").Value = '  +1.20%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.53%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.78%  '

# Row 18
$ws.Range("E18").Value = '  +0.67%  '

# Row 19
$ws.Range("E19").Value = '  -1.92%  '

# Row 20
$ws.Range("D20").Value = '71.223.73'
$ws.Range("E20").Value = '  +2.47%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '440.67'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.93%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.97%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '93.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.42%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.05%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.83%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.82%  '

# Row 28
$ws.Range("E28").Value = '  -0.83%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.39%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '690.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.52%  '

# Row 31
$ws.Range("E31").Value = '  +0.21%  '

# Row 32
$ws.Range("E32").Value = '  -0.12%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +13.98%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '67.68'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.14%  '

# Row 35
$ws.Range("D35").Value = '0.0₃0899'
$ws.Range("E35").Value = '  +1.22%  '

# Row 36
$ws.Range("E36").Value = '  -1.40%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.96'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.13%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.52'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.70%  '

# Row 39
$ws.Range("E39").Value = '  +0.14%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.27%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0494'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.59%  '

# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.13%  '

# Row 43
$ws.Range("E43").Value = '  -0.84%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.12'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.53%  '

# Row 45
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.10%  '

# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.23%  '

# Row 47
$ws.Range("E47").Value = '  +0.44%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000282'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +18.70%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.27'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.22%  '

# Row 50
$ws.Range("E50").Value = '  +0.14%  '

# Row 51
$ws.Range("E51").Value = '  -0.65%  '
